$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = 0.4906516671180725
$ws.Range("B1").Value = 0.3647835552692413
$ws.Range("C1").Value = 0.4184468984603882
$ws.Range("D1").Value = 3.903028249740601
$ws.Range("E1").Value = 1.635806679725647
